$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh COVID country stats + "datos actualizados" timestamp, and fix the handful of
# rows where the underlying data source reordered two neighbouring countries
# (Moldavia/Singapur, Jordania/Costa de Marfil/Republica de Macedonia, Albania/Eslovaquia,
# Jamaica/Eslovenia, Trinidad yTobago/Ruanda/Republica de Africa Central,
# Santa Lucia/Nueva Caledonia, Islas Malvinas/Montserrat).

# Footer timestamp (cell A1)
$ws.Range("A1").Value = 'Datos actualizados a 7 de Octubre de 2020 a las 17:30'

# Row 4
$ws.Range("B4").Value = 7730917
$ws.Range("C4").Value = 8171
$ws.Range("D4").Value = 4950124
$ws.Range("E4").Value = 2564729
$ws.Range("G4").Value = 242
$ws.Range("H4").Value = 216064

# Row 15
$ws.Range("B15").Value = 544275
$ws.Range("C15").Value = 14162
$ws.Range("G15").Value = 70
$ws.Range("H15").Value = 42515

# Row 17
$ws.Range("B17").Value = 474440
$ws.Range("C17").Value = 1134
$ws.Range("D17").Value = 447053
$ws.Range("E17").Value = 14297
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 13090

# Row 21
$ws.Range("B21").Value = 333940
$ws.Range("C21").Value = 3678
$ws.Range("D21").Value = 235303
$ws.Range("E21").Value = 62576
$ws.Range("G21").Value = 31
$ws.Range("H21").Value = 36061

# Row 26
$ws.Range("B26").Value = 309228
$ws.Range("C26").Value = 2109
$ws.Range("D26").Value = 267700
$ws.Range("E26").Value = 31883
$ws.Range("G26").Value = 10
$ws.Range("H26").Value = 9645

# Row 29
$ws.Range("B29").Value = 171906
$ws.Range("C29").Value = 583
$ws.Range("D29").Value = 144701
$ws.Range("E29").Value = 17674
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 9531

# Row 46
$ws.Range("B46").Value = 95704
$ws.Range("C46").Value = 834
$ws.Range("D46").Value = 84036
$ws.Range("E46").Value = 8333
$ws.Range("G46").Value = 25
$ws.Range("H46").Value = 3335

# Row 49
$ws.Range("B49").Value = 86543
$ws.Range("C49").Value = 496
$ws.Range("D49").Value = 79676
$ws.Range("E49").Value = 5262
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 1605

# Row 60: Moldavia
$ws.Range("A60").Value = 'Moldavia'
$ws.Range("B60").Value = 58794
$ws.Range("C60").Value = 1062
$ws.Range("D60").Value = 42480
$ws.Range("E60").Value = 14908
$ws.Range("G60").Value = 17
$ws.Range("H60").Value = 1406

# Row 61: Singapur
$ws.Range("A61").Value = 'Singapur'
$ws.Range("B61").Value = 57840
$ws.Range("C61").Value = 10
$ws.Range("D61").Value = 57624
$ws.Range("E61").Value = 189
$ws.Range("H61").Value = 27

# Row 88: Jordania
$ws.Range("A88").Value = 'Jordania'
$ws.Range("B88").Value = 20200
$ws.Range("C88").Value = 1199
$ws.Range("D88").Value = 5575
$ws.Range("E88").Value = 14494
$ws.Range("G88").Value = 9
$ws.Range("H88").Value = 131

# Row 89: Costa de Marfil
$ws.Range("A89").Value = 'Costa de Marfil'
$ws.Range("B89").Value = 19903
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 19539
$ws.Range("E89").Value = 244
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 120

# Row 90: Republica de Macedonia
$ws.Range("A90").Value = 'Republica de Macedonia'
$ws.Range("B90").Value = 19413
$ws.Range("C90").Value = 317
$ws.Range("D90").Value = 15749
$ws.Range("E90").Value = 2892
$ws.Range("G90").Value = 4
$ws.Range("H90").Value = 772

# Row 96: Albania
$ws.Range("A96").Value = 'Albania'
$ws.Range("B96").Value = 14730
$ws.Range("C96").Value = 162
$ws.Range("D96").Value = 9115
$ws.Range("E96").Value = 5208
$ws.Range("G96").Value = 4
$ws.Range("H96").Value = 407

# Row 97: Eslovaquia
$ws.Range("A97").Value = 'Eslovaquia'
$ws.Range("B97").Value = 14689
$ws.Range("C97").Value = 877
$ws.Range("D97").Value = 5200
$ws.Range("E97").Value = 9434
$ws.Range("H97").Value = 55

# Row 116: Jamaica
$ws.Range("A116").Value = 'Jamaica'
$ws.Range("B116").Value = 7191
$ws.Range("C116").Value = 82
$ws.Range("D116").Value = 2700
$ws.Range("E116").Value = 4365
$ws.Range("G116").Value = 3
$ws.Range("H116").Value = 126

# Row 117: Eslovenia
$ws.Range("A117").Value = 'Eslovenia'
$ws.Range("B117").Value = 7120
$ws.Range("C117").Value = 356
$ws.Range("D117").Value = 4535
$ws.Range("E117").Value = 2426
$ws.Range("H117").Value = 159

# Row 120
$ws.Range("B120").Value = 5898
$ws.Range("C120").Value = 15
$ws.Range("D120").Value = 5321
$ws.Range("E120").Value = 454

# Row 131: Trinidad yTobago
$ws.Range("A131").Value = 'Trinidad yTobago'
$ws.Range("B131").Value = 4876
$ws.Range("C131").Value = 30
$ws.Range("D131").Value = 3010
$ws.Range("E131").Value = 1782
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 84

# Row 132: Ruanda
$ws.Range("A132").Value = 'Ruanda'
$ws.Range("B132").Value = 4873
$ws.Range("D132").Value = 3246
$ws.Range("E132").Value = 1598
$ws.Range("H132").Value = 29

# Row 133: Republica de Africa Central
$ws.Range("A133").Value = 'Republica de Africa Central'
$ws.Range("B133").Value = 4852
$ws.Range("D133").Value = 1914
$ws.Range("E133").Value = 2876
$ws.Range("H133").Value = 62

# Row 207: Santa Lucia
$ws.Range("A207").Value = 'Santa Lucia'

# Row 208: Nueva Caledonia
$ws.Range("A208").Value = 'Nueva Caledonia'

# Row 215: Islas Malvinas
$ws.Range("A215").Value = 'Islas Malvinas'
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

# Row 216: Montserrat
$ws.Range("A216").Value = 'Montserrat'
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

Write-Host "Applied country/provincia data updates"
